$d = $word.ActiveDocument

# --- Change 1: expand the first paragraph with a red "(This is a change ...)" note ---
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.Text = "This is a Microsoft word document.  "

# Insert three separate red-colored runs after the existing text, before the
# paragraph mark, matching the original authoring (three runs split at the
# same boundaries as the source commit).
$chunk1 = [string][char]0x0028 + "This is a change " + [string][char]0x2013 + " Ve"
$chunk2 = "rsion for main branch"
$chunk3 = [string][char]0x0029

$pr = $d.Paragraphs.Item(1).Range
$ins = $d.Range($pr.Start, $pr.End - 1)
$ins.Collapse(0)
$ins.InsertAfter($chunk1)
$ins.Font.Color = 255

$pr2 = $d.Paragraphs.Item(1).Range
$ins2 = $d.Range($pr2.Start, $pr2.End - 1)
$ins2.Collapse(0)
$ins2.InsertAfter($chunk2)
$ins2.Font.Color = 255

$pr3 = $d.Paragraphs.Item(1).Range
$ins3 = $d.Range($pr3.Start, $pr3.End - 1)
$ins3.Collapse(0)
$ins3.InsertAfter($chunk3)
$ins3.Font.Color = 255

Write-Host "Paragraph 1 now:" $d.Paragraphs.Item(1).Range.Text

# --- Change 2: remove the trailing "ank God almighty, we are free at last." paragraph ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
Write-Host "Removing last paragraph:" $lastPara.Range.Text
$lastPara.Range.Delete()

Write-Host "New paragraph count:" $d.Paragraphs.Count
Write-Host "New last paragraph:" $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text
